$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (scraped data update): Price (column D) and
# Volume(1h) (column E) move for nearly every row, and two coin pairs
# (rows 13/14 and rows 30/32) swap rank order, so their Coin name (B),
# Link (C), Price (D) and Volume(1h) (E) cells all change together.
#
# Every write below is prefixed with a literal leading apostrophe. Excel
# treats that as the standard 'force text' marker, so number-looking
# strings already in the sheet (e.g. '226.61', '1.00', '34.402.59') stay
# literal text instead of being coerced into numeric cells -- matching
# how this data is stored (plain text) throughout columns B-E. The
# Style reset immediately after each write clears the transient quote-
# prefix flag so it doesn't linger as a spurious formatting change.

$r = $ws.Range('D2')
$r.Value = '''' + '34.402.59'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.Value = '''' + '  +0.86%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.Value = '''' + '1.795.98'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.Value = '''' + '  +0.55%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.Value = '''' + '  -0.01%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.Value = '''' + '226.61'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.Value = '''' + '  +0.17%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.Value = '''' + '0.554'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.Value = '''' + '  +1.37%  '
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.Value = '''' + '  -0.04%  '
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.Value = '''' + '32.41'
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.Value = '''' + '  +1.33%  '
$r.Style = "Normal"
$r = $ws.Range('D9')
$r.Value = '''' + '0.295'
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.Value = '''' + '  +1.35%  '
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.Value = '''' + '  +0.26%  '
$r.Style = "Normal"
$r = $ws.Range('D11')
$r.Value = '''' + '0.0951'
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.Value = '''' + '  +0.79%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.Value = '''' + '2.056.89'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.Value = '''' + '  +0.67%  '
$r.Style = "Normal"
$r = $ws.Range('B13')
$r.Value = '''' + 'WrappedEther'
$r.Style = "Normal"
$r = $ws.Range('C13')
$r.Value = '''' + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$r.Style = "Normal"
$r = $ws.Range('D13')
$r.Value = '''' + '1.822.27'
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.Value = '''' + '  +1.92%  '
$r.Style = "Normal"
$r = $ws.Range('B14')
$r.Value = '''' + 'Chainlink'
$r.Style = "Normal"
$r = $ws.Range('C14')
$r.Value = '''' + 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.Value = '''' + '11.08'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.Value = '''' + '  -1.15%  '
$r.Style = "Normal"
$r = $ws.Range('D15')
$r.Value = '''' + '0.629'
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.Value = '''' + '  +1.82%  '
$r.Style = "Normal"
$r = $ws.Range('D16')
$r.Value = '''' + '34.366.16'
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.Value = '''' + '  +0.90%  '
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.Value = '''' + '  +0.71%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.Value = '''' + '68.33'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.Value = '''' + '  +0.59%  '
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.Value = '''' + '  +3.14%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.Value = '''' + '246.51'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.Value = '''' + '  +0.40%  '
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.Value = '''' + '  +1.63%  '
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.Value = '''' + '  -0.07%  '
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.Value = '''' + '  +2.00%  '
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.Value = '''' + '2.05'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.Value = '''' + '  +0.28%  '
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.Value = '''' + '162.73'
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.Value = '''' + '  +0.67%  '
$r.Style = "Normal"
$r = $ws.Range('D26')
$r.Value = '''' + '7.19'
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.Value = '''' + '  +0.69%  '
$r.Style = "Normal"
$r = $ws.Range('D27')
$r.Value = '''' + '16.39'
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.Value = '''' + '  +0.49%  '
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.Value = '''' + '  +2.04%  '
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.Value = '''' + '  +0.12%  '
$r.Style = "Normal"
$r = $ws.Range('B30')
$r.Value = '''' + 'PancakeSwap'
$r.Style = "Normal"
$r = $ws.Range('C30')
$r.Value = '''' + 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$r.Style = "Normal"
$r = $ws.Range('D30')
$r.Value = '''' + '1.23'
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.Value = '''' + '  -0.07%  '
$r.Style = "Normal"
$r = $ws.Range('D31')
$r.Value = '''' + '0.0521'
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.Value = '''' + '  +0.81%  '
$r.Style = "Normal"
$r = $ws.Range('B32')
$r.Value = '''' + 'InternetComputer(DFINITY)'
$r.Style = "Normal"
$r = $ws.Range('C32')
$r.Value = '''' + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$r.Style = "Normal"
$r = $ws.Range('D32')
$r.Value = '''' + '3.89'
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.Value = '''' + '  +8.10%  '
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.Value = '''' + '  +3.45%  '
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.Value = '''' + '  +1.31%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.Value = '''' + '1.441.36'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.Value = '''' + '  -0.62%  '
$r.Style = "Normal"
$r = $ws.Range('D36')
$r.Value = '''' + '2.62'
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.Value = '''' + '  +9.32%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.Value = '''' + '0.664'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.Value = '''' + '  +3.05%  '
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.Value = '''' + '  +1.96%  '
$r.Style = "Normal"
$r = $ws.Range('D39')
$r.Value = '''' + '0.0191'
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.Value = '''' + '  -0.90%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.Value = '''' + '83.58'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.Value = '''' + '  +4.43%  '
$r.Style = "Normal"
$r = $ws.Range('D41')
$r.Value = '''' + '2.40'
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.Value = '''' + '  +1.30%  '
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.Value = '''' + '  +1.63%  '
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.Value = '''' + '  +2.88%  '
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.Value = '''' + '13.87'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.Value = '''' + '  +3.06%  '
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.Value = '''' + '0.0523'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.Value = '''' + '  +2.91%  '
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.Value = '''' + '  +0.91%  '
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.Value = '''' + '  +0.03%  '
$r.Style = "Normal"
$r = $ws.Range('D48')
$r.Value = '''' + '1.951.94'
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.Value = '''' + '  +0.38%  '
$r.Style = "Normal"
$r = $ws.Range('D49')
$r.Value = '''' + '105.73'
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.Value = '''' + '  -1.65%  '
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.Value = '''' + '  +0.02%  '
$r.Style = "Normal"
$r = $ws.Range('D51')
$r.Value = '''' + '0.0₆0126'
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.Value = '''' + '  -8.40%  '
$r.Style = "Normal"
